$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 18; Date = "2020-6-18"; Count = 0;     Label = $null },
    @{ Row = 19; Date = "2020-6-19"; Count = 0;     Label = $null },
    @{ Row = 20; Date = "2020-6-20"; Count = 20220; Label = "#COVIDIOTS: 20220" },
    @{ Row = 21; Date = "2020-6-21"; Count = 0;     Label = $null },
    @{ Row = 22; Date = "2020-6-23"; Count = 0;     Label = $null },
    @{ Row = 23; Date = "2020-6-24"; Count = 0;     Label = $null },
    @{ Row = 24; Date = "2020-6-25"; Count = 0;     Label = $null },
    @{ Row = 25; Date = "2020-6-26"; Count = 0;     Label = $null },
    @{ Row = 26; Date = "2020-6-27"; Count = 0;     Label = $null },
    @{ Row = 27; Date = "2020-6-28"; Count = 0;     Label = $null },
    @{ Row = 28; Date = "2020-6-29"; Count = 0;     Label = $null },
    @{ Row = 29; Date = "2020-6-30"; Count = 0;     Label = $null }
)

foreach ($item in $data) {
    $r = $item.Row
    $aCell = $ws.Cells.Item($r, 1)
    $aCell.Value = "'" + $item.Date
    $ws.Cells.Item($r, 2).Value = $item.Count
    if ($item.Label -ne $null) {
        $ws.Cells.Item($r, 3).Value = $item.Label
    }
}
